$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma-as-separator typos in provider name strings (first comma -> period) ---
$nameEdits = @(
    @('E22', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F22', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E72', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F72', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E103', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F103', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E130', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F130', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E183', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F183', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E24', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E73', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E71', 'FERNANDEZ MARIO H. GALLICET OSCAR M'),
    @('E89', 'FERNANDEZ MARIO H. GALLICET OSCAR M'),
    @('E92', 'TRABICHET MARIA. VERGARA ADEL Y OTRA'),
    @('F92', 'TRABICHET MARIA. VERGARA ADEL Y OTRA')
)
foreach ($edit in $nameEdits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

# --- Fix Importe amounts: "8.997,85" (es-AR) -> "8997.85" (plain decimal) ---
# These must stay TEXT cells (they were stored as text/shared-strings, not numbers),
# so force Text format before assigning, then drop back to the default style so no
# visible formatting changes - only the stored text for these cells changes.
$numEdits = @(
    @('H2', '8997.85'),
    @('H3', '15499.87'),
    @('H4', '19800.00'),
    @('H5', '348.00'),
    @('H6', '76.50'),
    @('H7', '2906.00'),
    @('H8', '59100.00'),
    @('H9', '122081.29'),
    @('H10', '5033.00'),
    @('H11', '1700.00'),
    @('H126', '1700.00'),
    @('H12', '8077.70'),
    @('H13', '3288.25'),
    @('H14', '9843.85'),
    @('H15', '749.50'),
    @('H16', '12926.77'),
    @('H17', '14715.00'),
    @('H18', '5.90'),
    @('H19', '178.40'),
    @('H20', '5726.99'),
    @('H21', '114.00'),
    @('H22', '22.20'),
    @('H23', '260.80'),
    @('H24', '250.00'),
    @('H25', '527.72'),
    @('H26', '5200.00'),
    @('H150', '5200.00'),
    @('H27', '10.60'),
    @('H28', '958.00'),
    @('H29', '443.00'),
    @('H30', '50190.00'),
    @('H31', '557.00'),
    @('H32', '9431.80'),
    @('H33', '3269.71'),
    @('H34', '1132.00'),
    @('H35', '1489.33'),
    @('H36', '724.33'),
    @('H37', '1359.08'),
    @('H38', '270.80'),
    @('H39', '3585.23'),
    @('H40', '13000.00'),
    @('H41', '87.02'),
    @('H42', '109.80'),
    @('H43', '47383.24'),
    @('H44', '1185.72'),
    @('H45', '4602.00'),
    @('H46', '23.50'),
    @('H47', '3715.92'),
    @('H48', '518.76'),
    @('H49', '598.04'),
    @('H50', '16500.00'),
    @('H51', '18900.00'),
    @('H52', '648.08'),
    @('H53', '300.00'),
    @('H54', '207.20'),
    @('H55', '707.40'),
    @('H56', '2521.17'),
    @('H57', '5814.24'),
    @('H58', '376.75'),
    @('H59', '80.00'),
    @('H60', '108.30'),
    @('H61', '1486.00'),
    @('H62', '2417.94'),
    @('H63', '3105.00'),
    @('H64', '180.00'),
    @('H131', '180.00'),
    @('H65', '650.00'),
    @('H66', '6240.00'),
    @('H67', '134.70'),
    @('H68', '152.50'),
    @('H69', '17148.00'),
    @('H70', '1692.00'),
    @('H71', '1651.00'),
    @('H72', '16197.00'),
    @('H73', '2909.80'),
    @('H74', '332.00'),
    @('H75', '1443.00'),
    @('H76', '241.32'),
    @('H77', '41.00'),
    @('H78', '86.20'),
    @('H79', '939.38'),
    @('H80', '350940.00'),
    @('H81', '808.00'),
    @('H82', '154.00'),
    @('H83', '117.00'),
    @('H84', '15.92'),
    @('H85', '7142.60'),
    @('H86', '3899.00'),
    @('H87', '6011.41'),
    @('H88', '4440.28'),
    @('H89', '12.50'),
    @('H90', '309.10'),
    @('H91', '17.00'),
    @('H92', '2725.00'),
    @('H93', '21.20'),
    @('H94', '3523.35'),
    @('H95', '1860.00'),
    @('H96', '1164.00'),
    @('H97', '54.30'),
    @('H98', '20.97'),
    @('H99', '4430.30'),
    @('H100', '48.97'),
    @('H101', '120.00'),
    @('H102', '40.00'),
    @('H103', '485.23'),
    @('H104', '1500.00'),
    @('H113', '1500.00'),
    @('H187', '1500.00'),
    @('H105', '10138.00'),
    @('H106', '1800.00'),
    @('H123', '1800.00'),
    @('H107', '5270.26'),
    @('H108', '364.00'),
    @('H109', '199.50'),
    @('H110', '116000.00'),
    @('H111', '4680.00'),
    @('H112', '700.00'),
    @('H114', '18827.20'),
    @('H115', '773.50'),
    @('H116', '500.00'),
    @('H121', '500.00'),
    @('H117', '600.00'),
    @('H118', '2600.00'),
    @('H119', '4000.00'),
    @('H120', '18453.48'),
    @('H122', '950.00'),
    @('H124', '1000.00'),
    @('H125', '1900.00'),
    @('H127', '200.00'),
    @('H128', '240.00'),
    @('H129', '9360.00'),
    @('H130', '258.38'),
    @('H132', '110.00'),
    @('H133', '4565.00'),
    @('H134', '38059.00'),
    @('H135', '6210.00'),
    @('H136', '6828.00'),
    @('H137', '595.36'),
    @('H138', '3041.00'),
    @('H139', '1421.00'),
    @('H140', '534.96'),
    @('H141', '279.00'),
    @('H142', '365.15'),
    @('H143', '17751.81'),
    @('H144', '8156.38'),
    @('H145', '4391.42'),
    @('H146', '50.26'),
    @('H147', '6329.05'),
    @('H148', '870.00'),
    @('H149', '999.16'),
    @('H151', '610.00'),
    @('H152', '357.70'),
    @('H153', '1940.09'),
    @('H154', '12.00'),
    @('H155', '440.00'),
    @('H156', '306.00'),
    @('H157', '603.00'),
    @('H158', '48814.96'),
    @('H159', '6142.09'),
    @('H160', '270.84'),
    @('H161', '3740.00'),
    @('H162', '2400.00'),
    @('H163', '2014.30'),
    @('H164', '70000.00'),
    @('H165', '508045.18'),
    @('H166', '4317.00'),
    @('H167', '16300.00'),
    @('H168', '1250.00'),
    @('H169', '96025.68'),
    @('H170', '215600.00'),
    @('H171', '17000.00'),
    @('H172', '247554.50'),
    @('H173', '10000.00'),
    @('H174', '232175.00'),
    @('H175', '223000.00'),
    @('H176', '230884.00'),
    @('H177', '226724.00'),
    @('H178', '193000.00'),
    @('H179', '104973.05'),
    @('H180', '77440.00'),
    @('H181', '29442.98'),
    @('H182', '9350.00'),
    @('H183', '25019.80'),
    @('H184', '12769.00'),
    @('H185', '55.00'),
    @('H186', '33500.00')
)
foreach ($edit in $numEdits) {
    $cell = $ws.Range($edit[0])
    $cell.NumberFormat = "@"
    $cell.Value = $edit[1]
    $cell.Style = "Normal"
}

